$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H8").Value = 6.5
$ws.Range("K8").Select()
